$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap Filecoin / FirstDigitalUSD rows (32 and 33) - name and link columns
$ws.Range("B32").Value = "Filecoin"
$ws.Range("B33").Value = "FirstDigitalUSD"

$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"

# Price column (D) updates; force text format for numeric-looking values to avoid auto-numeric conversion
$riskCells = @("D5","D6","D7","D9","D10","D11","D13","D15","D17","D19","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D41","D42","D44","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $riskCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.848.74"
$ws.Range("D3").Value = "2.255.29"
$ws.Range("D5").Value = "304.05"
$ws.Range("D6").Value = "94.94"
$ws.Range("D7").Value = "0.524"
$ws.Range("D9").Value = "0.487"
$ws.Range("D10").Value = "34.65"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D13").Value = "6.61"
$ws.Range("D14").Value = "2.603.98"
$ws.Range("D15").Value = "14.29"
$ws.Range("D16").Value = "2.257.81"
$ws.Range("D17").Value = "0.789"
$ws.Range("D18").Value = "41.741.64"
$ws.Range("D19").Value = "12.30"
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("D21").Value = "5.94"
$ws.Range("D22").Value = "67.88"
$ws.Range("D23").Value = "236.47"
$ws.Range("D24").Value = "2.56"
$ws.Range("D25").Value = "1.00"
$ws.Range("D26").Value = "1.92"
$ws.Range("D27").Value = "23.59"
$ws.Range("D28").Value = "36.31"
$ws.Range("D30").Value = "9.44"
$ws.Range("D31").Value = "159.27"
$ws.Range("D32").Value = "5.20"
$ws.Range("D33").Value = "0.999"
$ws.Range("D34").Value = "3.13"
$ws.Range("D35").Value = "0.0733"
$ws.Range("D36").Value = "17.07"
$ws.Range("D41").Value = "3.98"
$ws.Range("D42").Value = "2.38"
$ws.Range("D43").Value = "1.971.59"
$ws.Range("D44").Value = "0.0282"
$ws.Range("D46").Value = "2.91"
$ws.Range("D47").Value = "9.83"
$ws.Range("D48").Value = "53.02"
$ws.Range("D49").Value = "72.51"
$ws.Range("D50").Value = "1.50"
$ws.Range("D51").Value = "90.52"

# Volume(1h) column (E) updates
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +5.58%  "
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  -3.80%  "
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("E23").Value = "  -3.08%  "
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  +6.10%  "
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  -5.52%  "
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -0.96%  "
$ws.Range("E51").Value = "  -1.46%  "
